# Add new column 'event' to Card24 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Header: O1 = "event", formatted like the other header cells (N1) ---
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("O1").Value = "event"

# --- Row 2 previously had blank M2/N2 cells; they now carry values ---
$ws.Range("M2").Value = "تم تشغيل ماكينه"
$ws.Range("N2").Value = "لايوجد"

# --- New column O stays blank for every data row (2-12), but the cells
#     must still materialize in the sheet (matching the pre-existing blank
#     inline-string cells such as the old M2/N2). Assigning the (no-op)
#     "Normal" style forces Excel to keep an empty cell entry without
#     altering appearance or introducing a new style. ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Style = "Normal"
}
